$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's COM layer auto-detects plain ".Value" assignments that look like
# numbers (e.g. "1.00", "41.70") and silently coerces them to numeric cells,
# which drops significant trailing zeros and changes the stored cell type
# from text to number. The source data in this sheet is always text
# (t="inlineStr" originally), so force every write through a formula ->
# copy -> paste-values round trip: this evaluates to a genuine string result
# first, and pasting *values only* keeps it typed as Text without touching
# the cell's style/number format (avoids polluting styles.xml) and without
# leaving a literal leading apostrophe in the stored string.
function Set-TextValue($addr, $text) {
    $r = $ws.Range($addr)
    $escaped = $text.Replace('"', '""')
    $r.Formula = '="' + $escaped + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}

Set-TextValue "D2" "51.584.38"
Set-TextValue "E2" "  +4.85%  "
Set-TextValue "D3" "2.756.91"
Set-TextValue "E3" "  +4.94%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "116.42"
Set-TextValue "E5" "  +3.73%  "
Set-TextValue "D6" "332.32"
Set-TextValue "E6" "  +2.88%  "
Set-TextValue "D7" "0.539"
Set-TextValue "E7" "  +2.39%  "
Set-TextValue "E8" "  -0.07%  "
Set-TextValue "E9" "  +6.17%  "
Set-TextValue "D10" "41.70"
Set-TextValue "E10" "  +4.89%  "
Set-TextValue "D11" "0.0858"
Set-TextValue "E11" "  +5.83%  "
Set-TextValue "D12" "20.19"
Set-TextValue "E12" "  +2.32%  "
Set-TextValue "E13" "  +2.18%  "
Set-TextValue "D14" "7.65"
Set-TextValue "E14" "  +5.31%  "
Set-TextValue "D15" "3.188.61"
Set-TextValue "E15" "  +4.93%  "
Set-TextValue "D16" "2.755.41"
Set-TextValue "E16" "  +4.43%  "
Set-TextValue "D17" "0.886"
Set-TextValue "E17" "  +3.30%  "
Set-TextValue "D18" "51.534.76"
Set-TextValue "E18" "  +4.81%  "
Set-TextValue "E19" "  +6.44%  "
Set-TextValue "E20" "  +4.20%  "
Set-TextValue "E21" "  +2.54%  "
Set-TextValue "E22" "  +3.45%  "
Set-TextValue "D23" "278.19"
Set-TextValue "E24" "  +1.74%  "
Set-TextValue "E25" "  +4.61%  "
Set-TextValue "D26" "26.82"
Set-TextValue "E26" "  +2.54%  "
Set-TextValue "E28" "  -0.37%  "
Set-TextValue "D29" "2.22"
Set-TextValue "E29" "  -0.22%  "
Set-TextValue "D30" "0.141"
Set-TextValue "E30" "  +2.16%  "
Set-TextValue "D31" "35.13"
Set-TextValue "E31" "  +0.03%  "
Set-TextValue "E32" "  +0.86%  "
Set-TextValue "E33" "  +1.81%  "
Set-TextValue "E34" "  +2.76%  "
Set-TextValue "D35" "0.999"
Set-TextValue "E35" "  -0.21%  "
Set-TextValue "D36" "19.05"
Set-TextValue "E36" "  -0.06%  "
Set-TextValue "E37" "  +2.28%  "
Set-TextValue "D38" "4.98"
Set-TextValue "E38" "  +0.46%  "
Set-TextValue "D39" "3.23"
Set-TextValue "E39" "  +3.44%  "
Set-TextValue "E40" "  +10.56%  "
Set-TextValue "D41" "126.65"
Set-TextValue "E41" "  +0.58%  "
Set-TextValue "D42" "23.17"
Set-TextValue "E42" "  +4.49%  "
Set-TextValue "E43" "  +3.22%  "
Set-TextValue "D44" "2.30"
Set-TextValue "E44" "  +7.83%  "
Set-TextValue "D45" "2.43"
Set-TextValue "E45" "  +12.23%  "
Set-TextValue "D46" "2.090.10"
Set-TextValue "E46" "  +1.25%  "
Set-TextValue "E47" "  +3.29%  "
Set-TextValue "E48" "  +4.25%  "
Set-TextValue "E49" "  +6.45%  "
Set-TextValue "D50" "9.01"
Set-TextValue "E50" "  +1.13%  "
Set-TextValue "D51" "59.86"
Set-TextValue "E51" "  +1.64%  "

$excel.CutCopyMode = $false
